# Calculated average of variances instead of std's
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C.
# This shifts the old "Integral" column (C) to D, and old "Time" column (D) to E.
$ws.Columns("C:C").Insert()

# Header for the new Variance column
$ws.Range("C1").Value2 = "Variance"

# Variance formulas: Variance = STD^2 for every measurement row (2-11).
# Enter C2 on its own first (becomes a normal, non-shared formula),
# then fill C3:C11 together so Excel groups them into one shared formula.
$ws.Range("C2").Formula = "=B2^2"
$ws.Range("C3:C11").Formula = "=B3^2"

# Row 13 ("Avg"): add the average of the new Variance column
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"

# Row 14 ("STD"): no longer compute STD for the raw STD column (B) or the
# new Variance column (C) - only keep it for Integral (D) and Time (E)
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# New row 15: RMS = sqrt(average variance)
$ws.Range("A15").Value2 = "RMS"
$ws.Range("B15").Formula = "=SQRT(C13)"

# Match the new Variance column's width to column B's width
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth()

# Update the active selection like the author's saved view
$ws.Range("B16").Select() | Out-Null
